$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header / account holder details
$ws.Range("C2").Value = "Hartmut"

$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "2570314725427075"
$ws.Range("C3").Value = "Mohaupt"

# Opening balance line
$ws.Range("D5").Value = "KONTOSTAND AM 12.04.2024"

# Transaction rows
$ws.Range("B6").Value = "14.04."
$ws.Range("C6").Value = "15.04."
$ws.Range("D6").Value = "MITGLIEDSBEITRAG ZEUS BODYPOWER"
$ws.Range("E6").Value = "24,90-"

$ws.Range("B7").Value = "17.04."
$ws.Range("C7").Value = "18.04."
$ws.Range("D7").Value = "KARTENZAHLUNG ARAL TANKSTELLE"
$ws.Range("E7").Value = "58,42-"

$ws.Range("B8").Value = "19.04."
$ws.Range("C8").Value = "20.04."
$ws.Range("D8").Value = "RECHNUNG VODAFONE GMBH 30815065"
$ws.Range("E8").Value = "38,95-"

$ws.Range("B9").Value = "20.04."
$ws.Range("C9").Value = "21.04."
$ws.Range("D9").Value = "ABSCHLAG STROM Stadtwerke Rosenheim 87358683"
$ws.Range("E9").Value = "84,86-"

# Closing balance line
$ws.Range("D12").Value = "KONTOSTAND AM 24.04.2024"
$ws.Range("E12").Value = "207,13-"

# Next statement date
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 30.04.2024"
